# Append new scrape snapshot: 2026-01-19 06:35 JST
# Updates the "ランサーズ" sheet (案件情報) rows 2-8 with the newest data and
# adjusts a few column widths to fit the new content.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# --- Column width tweaks (raw OOXML widths: B=48, D=30, H=13) ---
# Excel's ColumnWidth property is offset from the stored sheet XML "width"
# by ~5/6 of a character (MDW padding), so subtract that to land exactly
# on the target stored width.
$ws.Columns.Item(2).ColumnWidth = 47.16666666666666   # stored width 48
$ws.Columns.Item(4).ColumnWidth = 29.16666666666667   # stored width 30
$ws.Columns.Item(8).ColumnWidth = 12.16666666666667   # stored width 13

# --- New timestamp applied to every refreshed row ---
$timestamp = "2026-01-19 06:35:35"

# --- Row 2: brand new listing ---
$ws.Range("A2").Value = $timestamp
$ws.Range("B2").Value = "大企業の業務効率化AIプロジェクトの技術方針策定を支援するAIテックリード募集"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("G2").Value = 385
$ws.Range("H2").Value = "🔥AI,Ai ◆効率化"

# --- Row 3 ---
$ws.Range("A3").Value = $timestamp
$ws.Range("B3").Value = "【医療機関向け業務改善サービスの新規開発】WEBアプリ開発におけるフルスタック開発担当者募集"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("G3").Value = 135
$ws.Range("H3").Value = "◆開発 ◇業務改善"

# --- Row 4 ---
$ws.Range("A4").Value = $timestamp
$ws.Range("B4").Value = "【急募】インバウンド向け新サービスアプリ開発見積作成依頼"
$ws.Range("C4").Value = "システム開発"
$ws.Range("D4").Value = "1,000 ~ 5,000 円 / 固定"
$ws.Range("E4").Value = "期限情報なし"
$ws.Range("G4").Value = 85
$ws.Range("H4").Value = "◆開発 ◇アプリ"

# --- Row 5 ---
$ws.Range("A5").Value = $timestamp
$ws.Range("B5").Value = "進行管理およびチームディレクションを担当"
$ws.Range("C5").Value = "システム開発"
$ws.Range("D5").Value = "~ 5,000 円 / 固定"
$ws.Range("E5").Value = "期限情報なし"
$ws.Range("G5").Value = 30
$ws.Range("H5").Value = "◇管理"

# --- Row 6 (loses its skill-summary column entirely) ---
$ws.Range("A6").Value = $timestamp
$ws.Range("B6").Value = "JSを使用したSaaSサービスの導入、保守のパートナー募集"
$ws.Range("C6").Value = "システム開発"
$ws.Range("D6").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E6").Value = "期限情報なし"
$ws.Range("G6").Value = 25
$ws.Range("H6").ClearContents()

# --- Row 7 ---
$ws.Range("A7").Value = $timestamp
$ws.Range("B7").Value = "【急募】自動車整備業向けCRM構築パートナー募集"
$ws.Range("C7").Value = "システム開発"
$ws.Range("D7").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("E7").Value = "期限情報なし"
$ws.Range("G7").Value = 25

# --- Row 8 ---
$ws.Range("A8").Value = $timestamp
$ws.Range("B8").Value = "移動型演出カートの電装設計および制御ユニット製作(Arduino/ESP32等)"
$ws.Range("C8").Value = "システム開発"
$ws.Range("D8").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("E8").Value = "期限情報なし"
$ws.Range("G8").Value = 18

# --- Refresh hyperlinks in column F (new target URLs) ---
$ws.Range("F2:F8").Hyperlinks.Delete()

$urls = @(
  "https://www.lancers.jp/work/detail/5423720",
  "https://www.lancers.jp/work/detail/5473940",
  "https://www.lancers.jp/work/detail/5474189",
  "https://www.lancers.jp/work/detail/5418064",
  "https://www.lancers.jp/work/detail/5474214",
  "https://www.lancers.jp/work/detail/5474125",
  "https://www.lancers.jp/work/detail/5474273"
)

for ($i = 0; $i -lt 7; $i++) {
  $row = 2 + $i
  $cell = $ws.Range("F$row")
  $ws.Hyperlinks.Add($cell, $urls[$i])
  $cell.Value = $urls[$i]
  $cell.Style = "Hyperlink"
}
